$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 960.60785
$ws.Cells.Item(15, 9).Value = 960.60785
$ws.Cells.Item(15, 11).Value = 2881.82355
$ws.Cells.Item(15, 13).Value = -2712.82355
$ws.Cells.Item(18, 8).Value = 2896
$ws.Cells.Item(18, 9).Value = 3133
$ws.Cells.Item(18, 10).Value = 1000
$ws.Cells.Item(18, 11).Value = 3133
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = -2849
$ws.Cells.Item(18, 14).Value = -1568
$ws.Cells.Item(19, 8).Value = 750.9474
$ws.Cells.Item(19, 10).Value = 1084
$ws.Cells.Item(19, 12).Value = 1084
$ws.Cells.Item(19, 14).Value = -1434
$ws.Cells.Item(88, 8).Value = 3125.9666
$ws.Cells.Item(88, 9).Value = 847
$ws.Cells.Item(88, 11).Value = 847
$ws.Cells.Item(88, 13).Value = -441
$ws.Cells.Item(91, 8).Value = 3125.9666
$ws.Cells.Item(91, 9).Value = 847
$ws.Cells.Item(91, 11).Value = 847
$ws.Cells.Item(91, 13).Value = 557
$ws.Cells.Item(92, 8).Value = 1375.6154
$ws.Cells.Item(92, 9).Value = 417.5238
$ws.Cells.Item(92, 11).Value = 417.5238
$ws.Cells.Item(92, 13).Value = 830.4762000000001
$ws.Cells.Item(94, 8).Value = 359
$ws.Cells.Item(94, 9).Value = 359
$ws.Cells.Item(94, 11).Value = 359
$ws.Cells.Item(94, 13).Value = 92
$ws.Cells.Item(96, 8).Value = 298.58334
$ws.Cells.Item(96, 10).Value = 204.25
$ws.Cells.Item(96, 12).Value = 612.75
$ws.Cells.Item(96, 14).Value = -3358.75
$ws.Cells.Item(99, 8).Value = 1256.091
$ws.Cells.Item(99, 10).Value = 1823.1428
$ws.Cells.Item(99, 12).Value = 5469.428400000001
$ws.Cells.Item(99, 14).Value = -8465.428400000001
$ws.Cells.Item(106, 8).Value = 2107.5
$ws.Cells.Item(106, 9).Value = 1649.25
$ws.Cells.Item(106, 10).Value = 3024
$ws.Cells.Item(106, 11).Value = 1649.25
$ws.Cells.Item(106, 12).Value = 3024
$ws.Cells.Item(106, 13).Value = -1018.25
$ws.Cells.Item(106, 14).Value = -4286
$ws.Cells.Item(132, 8).Value = 18520706
$ws.Cells.Item(132, 9).Value = 18870134
$ws.Cells.Item(132, 11).Value = 56610402
$ws.Cells.Item(132, 13).Value = -56607872
$ws.Cells.Item(137, 8).Value = 47535.31
$ws.Cells.Item(137, 9).Value = 58867.676
$ws.Cells.Item(137, 11).Value = 176603.028
$ws.Cells.Item(137, 13).Value = -174053.028
$ws.Cells.Item(138, 8).Value = 3026.2837
$ws.Cells.Item(138, 10).Value = 3546.9375
$ws.Cells.Item(138, 12).Value = 10640.8125
$ws.Cells.Item(138, 14).Value = -20920.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7763.3687
$ws.Cells.Item(32, 9).Value = 4393.2065
$ws.Cells.Item(32, 11).Value = 4393.2065
$ws.Cells.Item(32, 13).Value = -4106.2065
$ws.Cells.Item(74, 8).Value = 27394.527
$ws.Cells.Item(74, 9).Value = 1684.174
$ws.Cells.Item(74, 11).Value = 1684.174
$ws.Cells.Item(74, 13).Value = -810.174
$ws.Cells.Item(77, 8).Value = 27394.527
$ws.Cells.Item(77, 9).Value = 1684.174
$ws.Cells.Item(77, 11).Value = 8420.869999999999
$ws.Cells.Item(77, 13).Value = -4052.869999999999
$ws.Cells.Item(97, 8).Value = 1475957.8
$ws.Cells.Item(97, 9).Value = 2696850.5
$ws.Cells.Item(97, 10).Value = 10886.4
$ws.Cells.Item(97, 11).Value = 2696850.5
$ws.Cells.Item(97, 12).Value = 10886.4
$ws.Cells.Item(97, 13).Value = -2696354.5
$ws.Cells.Item(97, 14).Value = -11878.4
$ws.Cells.Item(102, 8).Value = 4905383.5
$ws.Cells.Item(102, 9).Value = 6412903.5
$ws.Cells.Item(102, 11).Value = 6412903.5
$ws.Cells.Item(102, 13).Value = -6411281.5
$ws.Cells.Item(110, 8).Value = 1324077
$ws.Cells.Item(110, 9).Value = 1324077
$ws.Cells.Item(110, 11).Value = 1324077
$ws.Cells.Item(110, 13).Value = -1322032
$ws.Cells.Item(132, 8).Value = 3290.9707
$ws.Cells.Item(132, 9).Value = 2520.4827
$ws.Cells.Item(132, 11).Value = 7561.4481
$ws.Cells.Item(132, 13).Value = -5031.4481

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 934.4
$ws.Cells.Item(11, 9).Value = 2500
$ws.Cells.Item(11, 10).Value = 543
$ws.Cells.Item(11, 11).Value = 2500
$ws.Cells.Item(11, 12).Value = 543
$ws.Cells.Item(11, 13).Value = -2360
$ws.Cells.Item(11, 14).Value = -823
$ws.Cells.Item(20, 8).Value = 8773792
$ws.Cells.Item(20, 9).Value = 15874520
$ws.Cells.Item(20, 10).Value = 2305.1765
$ws.Cells.Item(20, 11).Value = 15874520
$ws.Cells.Item(20, 12).Value = 2305.1765
$ws.Cells.Item(20, 13).Value = -15874273
$ws.Cells.Item(20, 14).Value = -2799.1765
$ws.Cells.Item(94, 8).Value = 4353032.5
$ws.Cells.Item(94, 9).Value = 6250719
$ws.Cells.Item(94, 10).Value = 15464
$ws.Cells.Item(94, 11).Value = 6250719
$ws.Cells.Item(94, 12).Value = 15464
$ws.Cells.Item(94, 13).Value = -6250268
$ws.Cells.Item(94, 14).Value = -16366
$ws.Cells.Item(107, 8).Value = 4204647
$ws.Cells.Item(107, 9).Value = 4764533.5
$ws.Cells.Item(107, 10).Value = 5500
$ws.Cells.Item(107, 11).Value = 4764533.5
$ws.Cells.Item(107, 12).Value = 5500
$ws.Cells.Item(107, 13).Value = -4762613.5
$ws.Cells.Item(107, 14).Value = -9340

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 32284.936
$ws.Cells.Item(31, 9).Value = 2501.9565
$ws.Cells.Item(31, 10).Value = 117911
$ws.Cells.Item(31, 11).Value = 2501.9565
$ws.Cells.Item(31, 12).Value = 117911
$ws.Cells.Item(31, 13).Value = -2206.9565
$ws.Cells.Item(31, 14).Value = -118501
$ws.Cells.Item(34, 8).Value = 32284.936
$ws.Cells.Item(34, 9).Value = 2501.9565
$ws.Cells.Item(34, 10).Value = 117911
$ws.Cells.Item(34, 11).Value = 2501.9565
$ws.Cells.Item(34, 12).Value = 117911
$ws.Cells.Item(34, 13).Value = -2299.9565
$ws.Cells.Item(34, 14).Value = -118315
$ws.Cells.Item(58, 8).Value = 7552.4287
$ws.Cells.Item(58, 9).Value = 11164
$ws.Cells.Item(58, 10).Value = 3579.7
$ws.Cells.Item(58, 11).Value = 11164
$ws.Cells.Item(58, 12).Value = 3579.7
$ws.Cells.Item(58, 13).Value = -10961
$ws.Cells.Item(58, 14).Value = -3985.7
$ws.Cells.Item(86, 8).Value = 14221.5
$ws.Cells.Item(86, 9).Value = 6777
$ws.Cells.Item(86, 10).Value = 21666
$ws.Cells.Item(86, 11).Value = 6777
$ws.Cells.Item(86, 12).Value = 21666
$ws.Cells.Item(86, 13).Value = -5654
$ws.Cells.Item(86, 14).Value = -23912
$ws.Cells.Item(89, 8).Value = 14221.5
$ws.Cells.Item(89, 9).Value = 6777
$ws.Cells.Item(89, 10).Value = 21666
$ws.Cells.Item(89, 11).Value = 33885
$ws.Cells.Item(89, 12).Value = 108330
$ws.Cells.Item(89, 13).Value = -28269
$ws.Cells.Item(89, 14).Value = -119562
$ws.Cells.Item(122, 8).Value = 2811.4443
$ws.Cells.Item(122, 9).Value = 2724.4375
$ws.Cells.Item(122, 11).Value = 8173.3125
$ws.Cells.Item(122, 13).Value = -5723.3125
$ws.Cells.Item(132, 8).Value = 65140.773
$ws.Cells.Item(132, 9).Value = 40869.48
$ws.Cells.Item(132, 11).Value = 122608.44
$ws.Cells.Item(132, 13).Value = -120078.44
$ws.Cells.Item(136, 8).Value = 7552.4287
$ws.Cells.Item(136, 9).Value = 11164
$ws.Cells.Item(136, 10).Value = 3579.7
$ws.Cells.Item(136, 11).Value = 33492
$ws.Cells.Item(136, 12).Value = 10739.1
$ws.Cells.Item(136, 13).Value = -30942
$ws.Cells.Item(136, 14).Value = -15839.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 77703.766
$ws.Cells.Item(5, 9).Value = 786.36365
$ws.Cells.Item(5, 10).Value = 500749.5
$ws.Cells.Item(5, 11).Value = 2359.09095
$ws.Cells.Item(5, 12).Value = 1502248.5
$ws.Cells.Item(5, 13).Value = -2247.09095
$ws.Cells.Item(5, 14).Value = -1502472.5
$ws.Cells.Item(8, 8).Value = 462.18182
$ws.Cells.Item(8, 9).Value = 462.18182
$ws.Cells.Item(8, 11).Value = 1386.54546
$ws.Cells.Item(8, 13).Value = -1247.54546
$ws.Cells.Item(37, 8).Value = 44533.332
$ws.Cells.Item(37, 10).Value = 44533.332
$ws.Cells.Item(37, 12).Value = 133599.996
$ws.Cells.Item(37, 14).Value = -133823.996
$ws.Cells.Item(56, 8).Value = 62504500
$ws.Cells.Item(56, 9).Value = 62504500
$ws.Cells.Item(56, 11).Value = 62504500
$ws.Cells.Item(56, 13).Value = -62503970
$ws.Cells.Item(122, 8).Value = 1359.625
$ws.Cells.Item(122, 10).Value = 2005
$ws.Cells.Item(122, 12).Value = 18045
$ws.Cells.Item(122, 14).Value = -22945
$ws.Cells.Item(131, 8).Value = 20837534
$ws.Cells.Item(131, 9).Value = 27784010
$ws.Cells.Item(131, 10).Value = 19611684
$ws.Cells.Item(131, 11).Value = 83352030
$ws.Cells.Item(131, 12).Value = 58835052
$ws.Cells.Item(131, 13).Value = -83346990
$ws.Cells.Item(131, 14).Value = -58845132
$ws.Cells.Item(135, 8).Value = 77703.766
$ws.Cells.Item(135, 9).Value = 786.36365
$ws.Cells.Item(135, 10).Value = 500749.5
$ws.Cells.Item(135, 11).Value = 7077.27285
$ws.Cells.Item(135, 12).Value = 4506745.5
$ws.Cells.Item(135, 13).Value = -4542.27285
$ws.Cells.Item(135, 14).Value = -4511815.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 5953506
$ws.Cells.Item(97, 9).Value = 11905512
$ws.Cells.Item(97, 10).Value = 1500
$ws.Cells.Item(97, 11).Value = 11905512
$ws.Cells.Item(97, 12).Value = 1500
$ws.Cells.Item(97, 13).Value = -11905016
$ws.Cells.Item(97, 14).Value = -2492
$ws.Cells.Item(102, 8).Value = 4386213
$ws.Cells.Item(102, 9).Value = 4832142.5
$ws.Cells.Item(102, 11).Value = 4832142.5
$ws.Cells.Item(102, 13).Value = -4830520.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1482.5
$ws.Cells.Item(16, 10).Value = 1999
$ws.Cells.Item(16, 12).Value = 1999
$ws.Cells.Item(16, 14).Value = -2339
$ws.Cells.Item(93, 8).Value = 13897746
$ws.Cells.Item(93, 9).Value = 16667200
$ws.Cells.Item(93, 10).Value = 50474.5
$ws.Cells.Item(93, 11).Value = 16667200
$ws.Cells.Item(93, 12).Value = 50474.5
$ws.Cells.Item(93, 13).Value = -16665952
$ws.Cells.Item(93, 14).Value = -52970.5
$ws.Cells.Item(100, 8).Value = 3933.1
$ws.Cells.Item(100, 9).Value = 3497.2856
$ws.Cells.Item(100, 11).Value = 3497.2856
$ws.Cells.Item(100, 13).Value = -2956.2856
$ws.Cells.Item(140, 8).Value = 98088.336
$ws.Cells.Item(140, 10).Value = 98087.5
$ws.Cells.Item(140, 12).Value = 98087.5
$ws.Cells.Item(140, 14).Value = -108447.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 4903.3794
$ws.Cells.Item(136, 9).Value = 4485.9
$ws.Cells.Item(136, 11).Value = 13457.7
$ws.Cells.Item(136, 13).Value = -10907.7
